# "Import user, items, warehouse, uom, department"
# Populates the first data row of the User import template (Sheet1) with a
# sample record and leaves the selection on the last cell touched (D3),
# matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: NIP, Name, Password, Role
$ws.Range("A3").Value = 123
$ws.Range("B3").Value = "gigi"
$ws.Range("C3").Value = "admin"
$ws.Range("D3").Value = 1

# Leave the active selection on D3, as in the saved workbook.
$ws.Range("D3").Select()
